# DEV 6 - CLI Changes for Manager Application and Registration
#
# ProjectApplication sheet:
#  - Row 3 (application for S1234567A) is updated to reflect the outcome of
#    the manager's decision: the application moves from "Pending" to
#    "Successful", the Application ID becomes 2, and the Date is refreshed
#    to the timestamp of the status change.
#  - The two blank spacer rows (4 and 5) are removed.
#  - The application that used to live on row 6 (T2345678D) shifts up to
#    row 4; its Application ID is renumbered from 1 to 3 now that it
#    follows the S1234567A application in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the S1234567A application row (row 3) ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("D3").Value = "Successful"
$ws.Range("F3").Value = 45770.013772708335

# --- Remove the blank placeholder rows ---
$ws.Rows("4:5").Delete()

# --- Renumber the application that shifted up into row 4 ---
$ws.Range("A4").Value = 3

# --- Restore the selection to where the CLI tool left the cursor ---
$ws.Range("F13").Select()
